# Reproduce the authored change:
#  - "Sheet1" (the data sheet with 3 rows) is renamed to "Data" and moved so it
#    becomes the first tab, ahead of "InputData". "UpdateInputData" stays last.
#  - That sheet keeps being the active/selected tab, but its selection moves
#    from A3 to E21.
$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$inputData = $wb.Worksheets.Item("InputData")

# Move "Sheet1" so it sits right before "InputData" (i.e. becomes the first tab).
$sheet1.Move($inputData) | Out-Null

# Re-fetch by name since the collection order changed.
$dataWs = $wb.Worksheets.Item("Sheet1")
$dataWs.Name = "Data"

# Keep it the active tab, but move the selection to E21.
$dataWs.Select() | Out-Null
$dataWs.Range("E21").Select() | Out-Null
